# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have now been handed back (in sync with en-US),
# filling in the "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns on the per-language sheets, and updating
# the status shown on the Overview sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": status columns for zh-cn (E) and de-de (F) move from
# "In Translation" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "zh-cn": fill in handback info for both files
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("I2").Value = "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md"
$wsZhCn.Range("J2").Value = "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.ca612527ec4c9233c0be0278e5400dc25a7272bc.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-27 04:24:12"

$wsZhCn.Range("I3").Value = "a6ac665e-6478-4485-ba74-ef18105598ee.md"
$wsZhCn.Range("J3").Value = "a6ac665e-6478-4485-ba74-ef18105598ee.55810bd326de1c3f50577f2c78476b5c873dc653.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-27 04:24:12"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/016bab939a51d6bdc8475cfa804ece3c3c1040bf/e2e/1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md", "", "", "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/016bab939a51d6bdc8475cfa804ece3c3c1040bf/e2e/a6ac665e-6478-4485-ba74-ef18105598ee.md", "", "", "a6ac665e-6478-4485-ba74-ef18105598ee.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": fill in handback info for both files
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("I2").Value = "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md"
$wsDeDe.Range("J2").Value = "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.ca612527ec4c9233c0be0278e5400dc25a7272bc.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-27 04:24:18"

$wsDeDe.Range("I3").Value = "a6ac665e-6478-4485-ba74-ef18105598ee.md"
$wsDeDe.Range("J3").Value = "a6ac665e-6478-4485-ba74-ef18105598ee.55810bd326de1c3f50577f2c78476b5c873dc653.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-27 04:24:18"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/016bab939a51d6bdc8475cfa804ece3c3c1040bf/e2e/1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md", "", "", "1fe8b380-449f-47e7-86cd-8cd1416cbdd4.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/016bab939a51d6bdc8475cfa804ece3c3c1040bf/e2e/a6ac665e-6478-4485-ba74-ef18105598ee.md", "", "", "a6ac665e-6478-4485-ba74-ef18105598ee.md") | Out-Null

# ---------------------------------------------------------------------
# Widen columns that now hold longer text so the new values are readable,
# mirroring the column auto-fit that happens when the report is regenerated.
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1
